# Scheduled-runner update: refreshes currentAveragePrice / Leve price /
# Leve profit figures (columns H-N) for a handful of Leve rows across the
# ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
#
# Columns: H=currentAveragePrice I=currentAveragePriceNQ J=currentAveragePriceHQ
#          K=LevePriceNQ L=LevePriceHQ M=LeveProfitNQ N=LeveProfitHQ
#
# Some rows don't carry every column (a blank LeveProfit column simply has
# no cell at all), so a few updates clear a stale cell or introduce a new
# one instead of just overwriting a value - ClearContents() / new Value
# assignments below reproduce that exactly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 69 - Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws.Range("H69").Value = 2103.25
$ws.Range("I69").Value = 2006.5
$ws.Range("J69").Value = 2200
$ws.Range("K69").Value = 6019.5
$ws.Range("L69").Value = 6600
$ws.Range("M69").Value = -5145.5
$ws.Range("N69").Value = -8348

# Row 72 - Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws.Range("H72").Value = 2103.25
$ws.Range("I72").Value = 2006.5
$ws.Range("J72").Value = 2200
$ws.Range("K72").Value = 18058.5
$ws.Range("L72").Value = 19800
$ws.Range("M72").Value = -13690.5
$ws.Range("N72").Value = -28536

# Row 132 - Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 558274.5600000001
$ws.Range("I132").Value = 3287.3076
$ws.Range("K132").Value = 9861.9228
$ws.Range("M132").Value = -7331.9228

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 21 - Fashion Weak / Iron Cuirass
$ws.Range("H21").Value = 13872.333
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 13872.333
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 13872.333
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -14620.333

# Row 32 - Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 7201.0205
$ws.Range("I32").Value = 6457.636
$ws.Range("J32").Value = 13742.8
$ws.Range("K32").Value = 6457.636
$ws.Range("L32").Value = 13742.8
$ws.Range("M32").Value = -6170.636
$ws.Range("N32").Value = -14316.8

# Row 110 - Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 471.33334
$ws.Range("I110").Value = 469.2857
$ws.Range("J110").Value = 500
$ws.Range("K110").Value = 469.2857
$ws.Range("L110").Value = 500
$ws.Range("M110").Value = 1575.7143
$ws.Range("N110").Value = -4590

# Row 122 - Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 1402.2354
$ws.Range("I122").Value = 1427.375
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4282.125
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -1832.125
$ws.Range("N122").Value = -7900

# Row 123 - The Armoire Is Open / High Durium Armguards of Maiming
$ws.Range("H123").Value = 47563.6
$ws.Range("J123").Value = 47563.6
$ws.Range("L123").Value = 47563.6
$ws.Range("N123").Value = -57363.6

# Row 132 - Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1994.2894
$ws.Range("I132").Value = 1438.6562
$ws.Range("K132").Value = 4315.9686
$ws.Range("M132").Value = -1785.9686

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 105 - Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 5145.7144
$ws.Range("I105").Value = 4624
$ws.Range("J105").Value = 6450
$ws.Range("K105").Value = 4624
$ws.Range("L105").Value = 6450
$ws.Range("M105").Value = -2877
$ws.Range("N105").Value = -9944

# ---------------------------------------------------------------- CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31 - Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 5919.381
$ws.Range("I31").Value = 1863.375
$ws.Range("J31").Value = 8415.385
$ws.Range("K31").Value = 1863.375
$ws.Range("L31").Value = 8415.385
$ws.Range("M31").Value = -1568.375
$ws.Range("N31").Value = -9005.385

# Row 34 - Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 5919.381
$ws.Range("I34").Value = 1863.375
$ws.Range("J34").Value = 8415.385
$ws.Range("K34").Value = 1863.375
$ws.Range("L34").Value = 8415.385
$ws.Range("M34").Value = -1661.375
$ws.Range("N34").Value = -8819.385

# ---------------------------------------------------------------- CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 62 - Little Orphan Candy / Fig Bavarois
$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Row 65 - Confections of Confession (L) / Fig Bavarois
$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# Row 82 - Persuasion of a Higher Power / Baked Pipira Pira
$ws.Range("H82").Value = 11312.5
$ws.Range("J82").Value = 12642.857
$ws.Range("L82").Value = 37928.571
$ws.Range("N82").Value = -38740.571

# Row 85 - Loaves and Fishes (L) / Baked Pipira Pira
$ws.Range("H85").Value = 11312.5
$ws.Range("J85").Value = 12642.857
$ws.Range("L85").Value = 37928.571
$ws.Range("N85").Value = -40736.571

# Row 113 - Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 582.34375
$ws.Range("I113").Value = 512.5
$ws.Range("J113").Value = 587
$ws.Range("K113").Value = 1537.5
$ws.Range("L113").Value = 1761
$ws.Range("M113").Value = 632.5
$ws.Range("N113").Value = -6101

# Row 122 - Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 2500522.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2500522.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 22504702.5
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -22509602.5

# ---------------------------------------------------------------- GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 6 - Bad Bromance / Bone Staff
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

# Row 16 - An Offer We Can't Refuse / Decorated Bone Staff
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

# Row 43 - Get the Green Stuff / Malachite Earrings
$ws.Range("H43").Value = 8727.375
$ws.Range("I43").Value = 6750
$ws.Range("J43").Value = 10704.75
$ws.Range("K43").Value = 6750
$ws.Range("L43").Value = 10704.75
$ws.Range("M43").Value = -6599
$ws.Range("N43").Value = -11006.75

# Row 130 - Planisphere to Paper / Chondrite Magitek Planisphere
$ws.Range("H130").Value = 45592
$ws.Range("J130").Value = 45592
$ws.Range("L130").Value = 45592
$ws.Range("N130").Value = -55632

# Row 131 - Star Athletes / Star Quartz Wristband of Aiming
$ws.Range("H131").Value = 25000
$ws.Range("J131").Value = 25000
$ws.Range("L131").Value = 25000
$ws.Range("N131").Value = -35080

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 7 - Tan Before the Ban / Leather
$ws.Range("H7").Value = 1293
$ws.Range("I7").Value = 1271.1818
$ws.Range("J7").Value = 1341
$ws.Range("K7").Value = 1271.1818
$ws.Range("L7").Value = 1341
$ws.Range("M7").Value = -1159.1818
$ws.Range("N7").Value = -1565

# Row 126 - Battered Books / Saiga Leather
$ws.Range("H126").Value = 1293
$ws.Range("I126").Value = 1271.1818
$ws.Range("J126").Value = 1341
$ws.Range("K126").Value = 3813.5454
$ws.Range("L126").Value = 4023
$ws.Range("M126").Value = -1343.5454
$ws.Range("N126").Value = -8963

# Row 134 - Freezing Fingers / Crocodileskin Fingerless Gloves of Striking
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# ---------------------------------------------------------------- WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 18 - Welcome to the Cotton Club / Cotton Halfgloves
$ws.Range("H18").Value = 806
$ws.Range("I18").Value = 806
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 806
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("M18").Value = -633

# Row 46 - Crunching the Numbers / Linen Hat
$ws.Range("H46").Value = 57441
$ws.Range("J46").Value = 57441
$ws.Range("L46").Value = 57441
$ws.Range("N46").Value = -57903

# Row 134 - Cloth for Canvas / Mountain Linen
$ws.Range("H134").Value = 57441
$ws.Range("J134").Value = 57441
$ws.Range("L134").Value = 172323
$ws.Range("N134").Value = -177393
